$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (2021-10-04, "Primera" quality) was added to the
# top of the data block. Insert a fresh row at row 5 (pushing the existing
# rows 5-28 down to 6-29) and populate it with the new entry's values.
$ws.Rows.Item(5).Insert()

$ws.Range("A5").Value = 5
$ws.Range("B5").Value = "Macroferia Regional de Talca"
$ws.Range("C5").Value = "Maule"
$ws.Range("D5").Value = 44473
$ws.Range("E5").Value = 7
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100107
$ws.Range("H5").Value = "Otros"
$ws.Range("I5").Value = 100107002
$ws.Range("J5").Value = "Chirimoya"
$ws.Range("K5").Value = "Cultivar IV Región"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 200
$ws.Range("N5").Value = 28000
$ws.Range("O5").Value = 28000
$ws.Range("P5").Value = 28000
$ws.Range("Q5").Value = "$/bandeja 10 kilos"
$ws.Range("R5").Value = "Provincia de Limarí"
$ws.Range("S5").Value = 2800
$ws.Range("T5").Value = 10
